# First stage in QC sprint 1
# The Student.xlsx template is being trimmed down: the Datebirth, Gender
# and Nationaid columns (C, D, E) are no longer collected, so remove them
# and let the remaining columns (Email, Password, Phone, Username) shift
# left into C:F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three obsolete header columns; Excel shifts everything after
# them (Email/Password/Phone/Username + their formatting) left automatically.
$ws.Range("C:E").Delete() | Out-Null

# Leave the selection where the author ended up after the edit.
$ws.Range("F8").Select() | Out-Null
